# Bump the "Förändrad" date (column C) by one day for all data rows (2-18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").End(4).Row  # xlDown -> last used row in column A
if ($lastRow -lt 18) { $lastRow = 18 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
